# Acceptance Test Plan update:
# Mark the "Listing Replays" and "Viewing Replays" acceptance criteria
# (rows 43-46 in the Test Plan sheet, Sprint 3 columns G/H) as Pass,
# with tester initials/date comment "SJD; 4/12/18".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

$rows = @(43, 44, 45, 46)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Pass"
    $ws.Cells.Item($r, 8).Value = "SJD; 4/12/18"
}
